# Fruta / hortaliza, semanal
# Insert a new weekly record as row 297 (pushing the existing rows 297-323
# down to 298-324) on the "Poroto granado" consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 297; everything below
# (previously 297..323) shifts down to 298..324.
$ws.Rows.Item(297).Insert()

# Populate the newly inserted row 297 with the new weekly observation.
$ws.Range("A297").Value = 9
$ws.Range("B297").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C297").Value = "Metropolitana"
$ws.Range("D297").Value = 44918
$ws.Range("E297").Value = 13
$ws.Range("F297").Value = 100112030
$ws.Range("G297").Value = "Poroto granado"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 70
$ws.Range("K297").Value = 36000
$ws.Range("L297").Value = 38000
$ws.Range("M297").Value = 37000
$ws.Range("N297").Value = '$/saco 25 kilos'
$ws.Range("O297").Value = "Región Metropolitana"
$ws.Range("P297").Value = 1480
$ws.Range("Q297").Value = 25
$ws.Range("R297").Value = "Hortaliza"
